$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# FirstName/LastName/Email/PhoneNumber/JoinDate columns one to the right.
$ws.Columns("A:A").Insert()

# Populate the new UserName column.
$ws.Range("A1").Value = "UserName"
$ws.Range("A2").Value = "admin"

# Match the original column width pattern (each original column kept its
# width, shifted one slot to the right); give the new column A a width.
$ws.Columns("A:A").ColumnWidth = 13.3

# The hyperlink that lived on the old Email column (now column D) needs to
# be re-anchored; the engine does not auto-shift hyperlink anchors when
# columns are inserted, so recreate it on the correct cell.
$null = $ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("D2"), "mailto:admin@admin.com")
$ws.Range("D2").Style = "Hyperlink"

# Restore the selection to match the saved view state.
$null = $ws.Range("B12").Select()
